$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C3"
$ws.Range("C2").Value = "C3ar1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 46.85851166666667
$ws.Range("H2").Value = 140.575535
$ws.Range("I2").Value = 0.1419057303676978
$ws.Range("J2").Value = 0.1419057303676978
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.91561933333333
$ws.Range("N2").Value = 44.746858
$ws.Range("O2").Value = 0.9776685906860466
$ws.Range("P2").Value = 0.9776685906860467
$ws.Range("Q2").Value = 698.9237225465588
$ws.Range("R2").Value = 6290.31350291903
$ws.Range("S2").Value = 0.1387367754188612
$ws.Range("T2").Value = 0.1387367754188613

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C3"
$ws.Range("C3").Value = "C3ar1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 46.85851166666667
$ws.Range("H3").Value = 140.575535
$ws.Range("I3").Value = 0.1419057303676978
$ws.Range("J3").Value = 0.1419057303676978
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.340695
$ws.Range("N3").Value = 1.022085
$ws.Range("O3").Value = 0.02233140931395335
$ws.Range("P3").Value = 0.02233140931395336
$ws.Range("Q3").Value = 15.964460632275
$ws.Range("R3").Value = 143.680145690475
$ws.Range("S3").Value = 0.00316895494883656
$ws.Range("T3").Value = 0.003168954948836561

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "C3"
$ws.Range("C4").Value = "C3ar1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 283.1772156666667
$ws.Range("H4").Value = 849.531647
$ws.Range("I4").Value = 0.8575703363889615
$ws.Range("J4").Value = 0.8575703363889616
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.91561933333333
$ws.Range("N4").Value = 44.746858
$ws.Range("O4").Value = 0.9776685906860466
$ws.Range("P4").Value = 0.9776685906860467
$ws.Range("Q4").Value = 4223.763552757236
$ws.Range("R4").Value = 38013.87197481513
$ws.Range("S4").Value = 0.8384195821915549
$ws.Range("T4").Value = 0.8384195821915551

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "C3"
$ws.Range("C5").Value = "C3ar1"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 283.1772156666667
$ws.Range("H5").Value = 849.531647
$ws.Range("I5").Value = 0.8575703363889615
$ws.Range("J5").Value = 0.8575703363889616
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.340695
$ws.Range("N5").Value = 1.022085
$ws.Range("O5").Value = 0.02233140931395335
$ws.Range("P5").Value = 0.02233140931395336
$ws.Range("Q5").Value = 96.477061491555
$ws.Range("R5").Value = 868.2935534239951
$ws.Range("S5").Value = 0.01915075419740657
$ws.Range("T5").Value = 0.01915075419740657

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "C3"
$ws.Range("C6").Value = "C3ar1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1730073333333333
$ws.Range("H6").Value = 0.519022
$ws.Range("I6").Value = 0.000523933243340694
$ws.Range("J6").Value = 0.000523933243340694
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.91561933333333
$ws.Range("N6").Value = 44.746858
$ws.Range("O6").Value = 0.9776685906860466
$ws.Range("P6").Value = 0.9776685906860467
$ws.Range("Q6").Value = 2.580511525875111
$ws.Range("R6").Value = 23.224603732876
$ws.Range("S6").Value = 0.0005122330756304658
$ws.Range("T6").Value = 0.0005122330756304659

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "C3"
$ws.Range("C7").Value = "C3ar1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1730073333333333
$ws.Range("H7").Value = 0.519022
$ws.Range("I7").Value = 0.000523933243340694
$ws.Range("J7").Value = 0.000523933243340694
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.340695
$ws.Range("N7").Value = 1.022085
$ws.Range("O7").Value = 0.02233140931395335
$ws.Range("P7").Value = 0.02233140931395336
$ws.Range("Q7").Value = 0.05894273343
$ws.Range("R7").Value = 0.53048460087
$ws.Range("S7").Value = 0.0000117001677102281595330363420104369254
$ws.Range("T7").Value = 0.0000117001677102281595330363420104369254

